# Add the newly evaluated "reduced test size" Z3alpha stats for
# QF_NIA, UFNIA and QF_SLIA parallel-strategy evaluations.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# UFNIA row (row 13): Z3alpha column F
$ws.Range("F13").Value = "2483 [4core, 60s timeout, 5000 instances]"

# QF_NIA row (row 20): Z3alpha column F
$ws.Range("F20").Value = "3968 [4core, 60s timeout, 5000 instances]"

# QF_SLIA row (row 24): Z3alpha column F
$ws.Range("F24").Value = "9418 [4core, 60s timeout, 10000 instances]"

# Leave the cursor on the last-edited cell, matching the author's saved view state.
$ws.Range("F24").Select()
